# Generate Report for Handback
# Refreshes the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) for the first data row
# (row 2) of each language sheet with new report timestamps.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-22 18:49:57"
$ws_zhcn.Range("H2").Value = "2016-03-22 18:50:30"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-22 18:50:01"
$ws_dede.Range("H2").Value = "2016-03-22 18:50:38"
